$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last refreshed" timestamp caption (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 00:22"

# --- Update country statistics (columns B:H = Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 500879
$ws.Range("C4").Value = 32313
$ws.Range("D4").Value = 27239
$ws.Range("E4").Value = 455003
$ws.Range("F4").Value = 10916
$ws.Range("G4").Value = 1946
$ws.Range("H4").Value = 18637

# Row 5 - Espana
$ws.Range("B5").Value = 158273
$ws.Range("C5").Value = 5051
$ws.Range("D5").Value = 55668
$ws.Range("E5").Value = 86524
$ws.Range("F5").Value = 7371
$ws.Range("G5").Value = 634
$ws.Range("H5").Value = 16081

# Row 8 - Alemania
$ws.Range("B8").Value = 122171
$ws.Range("C8").Value = 3936
$ws.Range("D8").Value = 53913
$ws.Range("E8").Value = 65491
$ws.Range("F8").Value = 4895
$ws.Range("G8").Value = 160
$ws.Range("H8").Value = 2767

# Row 9 - China
$ws.Range("B9").Value = 81907
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 77455
$ws.Range("E9").Value = 1116
$ws.Range("F9").Value = 144
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 3336

# Row 10 - Reino Unido
$ws.Range("B10").Value = 73758
$ws.Range("C10").Value = 8681
$ws.Range("D10").Value = 344
$ws.Range("E10").Value = 64456
$ws.Range("F10").Value = 1559
$ws.Range("G10").Value = 980
$ws.Range("H10").Value = 8958

# Rows 105/106 - Nigeria and Kirguistan swap positions; Nigeria (now row 105)
# receives updated totals, Kirguistan (now row 106) keeps its previous totals.
$ws.Range("A105").Value = "Nigeria"
$ws.Range("B105").Value = 305
$ws.Range("C105").Value = 17
$ws.Range("D105").Value = 58
$ws.Range("E105").Value = 240
$ws.Range("F105").Value = 2
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 7

$ws.Range("A106").Value = "Kirguistan"
$ws.Range("B106").Value = 298
$ws.Range("C106").Value = 18
$ws.Range("D106").Value = 35
$ws.Range("E106").Value = 258
$ws.Range("F106").Value = 5
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 5

# Row 121 - Martinica
$ws.Range("B121").Value = 155
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 50
$ws.Range("E121").Value = 99
$ws.Range("F121").Value = 19
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 6
